# Workbook/worksheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Renumber the existing "index" column (A) for rows 46-52 -----------
# These rows had a gap in the sequential counter (45,47,48,49,50,51,52) that
# gets closed up (44,45,46,47,48,49,50) now that extra rows are appended
# below with their own (also slightly gappy) numbering.
$ws.Range("A46").Value = 44
$ws.Range("A47").Value = 45
$ws.Range("A48").Value = 46
$ws.Range("A49").Value = 47
$ws.Range("A50").Value = 48
$ws.Range("A51").Value = 49
$ws.Range("A52").Value = 50

# --- 2) Copy the "index" column formatting (bold/border/center style) down
# onto the new rows 53-59 so they match the rest of column A, without
# disturbing the values we are about to write explicitly. -----------------
$ws.Range("A52").Copy()
$ws.Range("A53:A59").PasteSpecial(-4122)

# --- 3) Append the new rows 53-59 with their data --------------------------
$ws.Range("A53").Value = 52
$ws.Range("B53").Value = 1.775319622012229
$ws.Range("C53").Value = 3193.8
$ws.Range("D53").Value = 0.01634241245136187
$ws.Range("E53").Value = 29.4
$ws.Range("F53").Value = 41
$ws.Range("G53").Value = "Portgas D Åce "
$ws.Range("H53").Value = "SOLO"
$ws.Range("I53").Value = 0.02279043913285158
$ws.Range("J53").Value = 4.4
$ws.Range("K53").Value = 0.002445803224013341

$ws.Range("A54").Value = 53
$ws.Range("B54").Value = 7.410958116892823
$ws.Range("C54").Value = 13175.2
$ws.Range("D54").Value = 0.05164043510720542
$ws.Range("E54").Value = 91.2
$ws.Range("F54").Value = 278.8
$ws.Range("G54").Value = "LS DUFFY"
$ws.Range("H54").Value = "SOLO"
$ws.Range("I54").Value = 0.1567269453783251
$ws.Range("J54").Value = 8.199999999999999
$ws.Range("K54").Value = 0.004627791059971394

$ws.Range("A55").Value = 54
$ws.Range("B55").Value = 6.474714026753915
$ws.Range("C55").Value = 8602
$ws.Range("D55").Value = 0.04634044514716931
$ws.Range("E55").Value = 61.2
$ws.Range("F55").Value = 328.6
$ws.Range("G55").Value = "BigFather Rengar"
$ws.Range("H55").Value = "SOLO"
$ws.Range("I55").Value = 0.2473481266503256
$ws.Range("J55").Value = 9.6
$ws.Range("K55").Value = 0.007220967752454279

$ws.Range("A56").Value = 56
$ws.Range("B56").Value = 1.901710291787398
$ws.Range("C56").Value = 4185.8
$ws.Range("D56").Value = 0.01349240180076153
$ws.Range("E56").Value = 29.6
$ws.Range("F56").Value = 204.4
$ws.Range("G56").Value = "MyDogaN"
$ws.Range("H56").Value = "DUO_SUPPORT"
$ws.Range("I56").Value = 0.09128055152368671
$ws.Range("J56").Value = 18.6
$ws.Range("K56").Value = 0.008167217339014521

$ws.Range("A57").Value = 57
$ws.Range("B57").Value = 2.791646791513082
$ws.Range("C57").Value = 3890.6
$ws.Range("D57").Value = 0.01949175076877424
$ws.Range("E57").Value = 28.2
$ws.Range("F57").Value = 221.4
$ws.Range("G57").Value = "Mr Kayn"
$ws.Range("H57").Value = "DUO_SUPPORT"
$ws.Range("I57").Value = 0.1749684198889241
$ws.Range("J57").Value = 3.2
$ws.Range("K57").Value = 0.002463335941977546

$ws.Range("A58").Value = 58
$ws.Range("B58").Value = 0
$ws.Range("C58").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = ""
$ws.Range("H58").Value = "SOLO"
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0

$ws.Range("A59").Value = 59
$ws.Range("B59").Value = 1.775319622012229
$ws.Range("C59").Value = 3193.8
$ws.Range("D59").Value = 0.01634241245136187
$ws.Range("E59").Value = 29.4
$ws.Range("F59").Value = 41
$ws.Range("G59").Value = "Portgas D Åce "
$ws.Range("H59").Value = "SOLO"
$ws.Range("I59").Value = 0.02279043913285158
$ws.Range("J59").Value = 4.4
$ws.Range("K59").Value = 0.002445803224013341
